$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.50609999999998
$ws.Range("E6").Value = 16.56789999999999
$ws.Range("A7").Value = -21.9068
$ws.Range("B7").Value = 4.971200000000001
$ws.Range("B15").Value = 4.738899999999994
$ws.Range("A16").Value = -21.9113
$ws.Range("C16").Value = -12.639
$ws.Range("C19").Value = -13.02980000000001
$ws.Range("E19").Value = 16.31609999999999
$ws.Range("E20").Value = 16.1478
$ws.Range("B21").Value = 10.6061
$ws.Range("E21").Value = 16.39650000000002
$ws.Range("B22").Value = 10.0795
$ws.Range("B23").Value = 8.776999999999999
$ws.Range("E24").Value = 16.18069999999999
$ws.Range("A28").Value = -21.9659
$ws.Range("A29").Value = -21.57169999999998
$ws.Range("A32").Value = -21.219
$ws.Range("B34").Value = 8.302000000000007
$ws.Range("D34").Value = -8.084500000000002
$ws.Range("E35").Value = 16.09139999999999
$ws.Range("C36").Value = -12.75720000000001
$ws.Range("E39").Value = 15.42609999999999
$ws.Range("A40").Value = -19.33899999999999
$ws.Range("E41").Value = 16.51599999999998
$ws.Range("B43").Value = 5.225399999999999
$ws.Range("D43").Value = -8.303700000000003
$ws.Range("B45").Value = 5.1849
$ws.Range("C46").Value = -14.61519999999999
$ws.Range("D48").Value = -7.395999999999997
$ws.Range("B50").Value = 4.843899999999997
$ws.Range("C50").Value = -13.58699999999999
$ws.Range("B51").Value = 5.802899999999999
$ws.Range("A52").Value = -22.285
$ws.Range("A57").Value = -22.2694
$ws.Range("A66").Value = -21.4781
$ws.Range("B66").Value = 5.666999999999998
$ws.Range("B67").Value = 5.259299999999999
$ws.Range("D70").Value = -7.441799999999997
$ws.Range("D73").Value = -7.901099999999999
$ws.Range("E73").Value = 16.34069999999999
$ws.Range("E76").Value = 16.24699999999999
$ws.Range("B79").Value = 9.773700000000005
$ws.Range("B84").Value = 5.1897
$ws.Range("E85").Value = 16.3587
$ws.Range("D87").Value = -8.499799999999995
$ws.Range("B92").Value = 4.665599999999998
$ws.Range("D92").Value = -6.215800000000001
$ws.Range("C95").Value = -12.2179
$ws.Range("B97").Value = 6.376699999999995
$ws.Range("C97").Value = -12.4018
$ws.Range("E98").Value = 16.3665
$ws.Range("A100").Value = -22.1847
$ws.Range("D101").Value = -8.151
